$wb = $excel.ActiveWorkbook

$wsEmployees = $wb.Worksheets.Item("Employees")
$wsEmployees.Range("B9").Value = "E0127"

$wsRoles = $wb.Worksheets.Item("Roles")
$wsRoles.Range("B3").ClearContents()
$wsRoles.Range("B5").Value = "Data engineer 4"
$wsRoles.Range("B6").ClearContents()
